$wb = $excel.ActiveWorkbook

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "GNG_TO-16509961913490465"
$ws1.Range("B2").Value = "go_stims-16509961913170722.csv"
$ws1.Range("B3").Value = "GNG_stims-16509961913330755.csv"
$ws1.Range("B4").Value = "go_stims-16509961913330755.csv"
$ws1.Range("B5").Value = "GNG_stims-16509961913490465.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16509961940931132"
$ws2.Range("B2").Value = "TB-16509961931650734.csv"
$ws2.Range("B3").Value = "ZB-match_6-1650996192109042.csv"
$ws2.Range("B4").Value = "ZB-match_7-1650996192173087.csv"
$ws2.Range("B5").Value = "OB-16509961924930778.csv"
$ws2.Range("B6").Value = "TB-1650996194069082.csv"
$ws2.Range("B7").Value = "ZB-match_6-16509961914130428.csv"
$ws2.Range("B8").Value = "TB-16509961939330823.csv"
$ws2.Range("B9").Value = "OB-16509961926930804.csv"
$ws2.Range("B10").Value = "OB-16509961922850437.csv"

# --- Sheet 3: RS_TO (name only) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16509961940931132"

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16509961941410816"
$ws4.Range("B2").Value = "MM_stims-16509961941090858.csv"
$ws4.Range("B3").Value = "ZM_stims-16509961940931132.csv"
$ws4.Range("B4").Value = "MM_stims-16509961941250787.csv"
$ws4.Range("B5").Value = "ZM_stims-16509961941090858.csv"
$ws4.Range("B6").Value = "MM_stims-16509961941410816.csv"
$ws4.Range("B7").Value = "ZM_stims-16509961941250787.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "vSAT_TO-1650996194221098"
$ws5.Range("B2").Value = "SAT_stims-165099619414905.csv"
$ws5.Range("B3").Value = "vSAT_stims-16509961941890543.csv"
$ws5.Range("B4").Value = "SAT_stims-16509961941730843.csv"
$ws5.Range("B5").Value = "vSAT_stims-16509961942050872.csv"
